$wb = $excel.ActiveWorkbook

$wsProd = $wb.Worksheets.Item("Productdata")
$wsCap = $wb.Worksheets.Item("Capacity")
$wsProc = $wb.Worksheets.Item("ProcessingTime")

# Productdata
$wsProd.Range("C2").Value = 0
$wsProd.Range("E2").Value = 69.61767749999999
$wsProd.Range("C3").Value = 0
$wsProd.Range("E3").Value = 51.421311
$wsProd.Range("C4").Value = 0
$wsProd.Range("E4").Value = 10.543545
$wsProd.Range("C5").Value = 0
$wsProd.Range("E5").Value = 104.7326625
$wsProd.Range("C6").Value = 0
$wsProd.Range("E6").Value = 11.227599
$wsProd.Range("C7").Value = 0
$wsProd.Range("E7").Value = 91.39680000000001
$wsProd.Range("C8").Value = 0
$wsProd.Range("E8").Value = 5.893560000000001
$wsProd.Range("C9").Value = 0
$wsProd.Range("E9").Value = 9.135017999999999
$wsProd.Range("C10").Value = 0
$wsProd.Range("E10").Value = 8.840339999999999
$wsProd.Range("C11").Value = 0
$wsProd.Range("E11").Value = 6.261907500000001
$wsProd.Range("C12").Value = 0
$wsProd.Range("E12").Value = 6.114568499999999
$wsProd.Range("C13").Value = 0
$wsProd.Range("E13").Value = 8.7666705
$wsProd.Range("C14").Value = 0
$wsProd.Range("E14").Value = 0.3870225
$wsProd.Range("C15").Value = 0
$wsProd.Range("E15").Value = 58.14009900000001
$wsProd.Range("C16").Value = 0
$wsProd.Range("E16").Value = 4.0383405
$wsProd.Range("C17").Value = 0
$wsProd.Range("E17").Value = 6.114568499999999
$wsProd.Range("C18").Value = 0
$wsProd.Range("E18").Value = 0.458343
$wsProd.Range("C19").Value = 0
$wsProd.Range("E19").Value = 57.8103885
$wsProd.Range("C20").Value = 0
$wsProd.Range("E20").Value = 2.513412
$wsProd.Range("C21").Value = 2580
$wsProd.Range("E21").Value = 4.6463625
$wsProd.Range("C22").Value = 3177
$wsProd.Range("E22").Value = 6.400718999999999
$wsProd.Range("C23").Value = 21144
$wsProd.Range("E23").Value = 49.49437499999999
$wsProd.Range("C24").Value = 25613
$wsProd.Range("E24").Value = 62.9124975
$wsProd.Range("C25").Value = 3068
$wsProd.Range("E25").Value = 6.371662500000001
$wsProd.Range("C26").Value = 2514
$wsProd.Range("E26").Value = 5.344254
$wsProd.Range("C27").Value = 16254
$wsProd.Range("E27").Value = 42.21759150000001
$wsProd.Range("C28").Value = 22320
$wsProd.Range("E28").Value = 54.34421399999999

# Capacity
$wsCap.Range("B2").Value = 327420
$wsCap.Range("B3").Value = 818550
$wsCap.Range("B4").Value = 16700
$wsCap.Range("B5").Value = 83870
$wsCap.Range("B6").Value = 40100
$wsCap.Range("B7").Value = 253880
$wsCap.Range("B8").Value = 327420
$wsCap.Range("B9").Value = 654840
$wsCap.Range("B10").Value = 491130
$wsCap.Range("B11").Value = 818550
$wsCap.Range("B12").Value = 163710
$wsCap.Range("B13").Value = 491130
$wsCap.Range("B14").Value = 8350
$wsCap.Range("B15").Value = 461100
$wsCap.Range("B16").Value = 83870
$wsCap.Range("B17").Value = 654840
$wsCap.Range("B18").Value = 24060
$wsCap.Range("B19").Value = 142980
$wsCap.Range("B20").Value = 190410
$wsCap.Range("B21").Value = 17650
$wsCap.Range("B22").Value = 4820
$wsCap.Range("B23").Value = 112500
$wsCap.Range("B24").Value = 185480
$wsCap.Range("B25").Value = 8700
$wsCap.Range("B26").Value = 3670
$wsCap.Range("B27").Value = 111720
$wsCap.Range("B28").Value = 35540

# ProcessingTime
$wsProc.Range("B2").Value = 2
$wsProc.Range("C3").Value = 5
$wsProc.Range("E5").Value = 1
$wsProc.Range("F6").Value = 5
$wsProc.Range("G7").Value = 4
$wsProc.Range("H8").Value = 2
$wsProc.Range("I9").Value = 4
$wsProc.Range("J10").Value = 3
$wsProc.Range("K11").Value = 5
$wsProc.Range("M13").Value = 3
$wsProc.Range("O15").Value = 5
$wsProc.Range("P16").Value = 1
$wsProc.Range("Q17").Value = 4
$wsProc.Range("S19").Value = 2
$wsProc.Range("V22").Value = 1
$wsProc.Range("W23").Value = 3
$wsProc.Range("Z26").Value = 1
$wsProc.Range("AA27").Value = 4
$wsProc.Range("AB28").Value = 1
